$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 2025-02-10 to 2025-02-11 (serial 45698 -> 45699)
# for every data row (2-37).
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}

# Remove the last data row (row 38) entirely.
$ws.Rows("38:38").Delete()

# Row 37 (now the last row) loses its explicit custom-height flag, matching
# the row that used to be the trailing row before this edit.
$ws.Rows("37:37").AutoFit()
